$wb = $excel.ActiveWorkbook

# Locate the "Portugal" sheet, which is the template for the new "Slovakia" sheet.
$portugal = $wb.Worksheets.Item("Portugal")

# Copy Portugal to create the new sheet, placing it right after Portugal.
$portugal.Copy($null, $portugal)

# The newly created sheet becomes the active sheet, placed right after Portugal.
$slovakia = $wb.ActiveSheet
$slovakia.Name = "Slovakia"

# Update the market name and user story cells for the new Slovakia sheet.
$slovakia.Range("B4").Value = "NGC-2930/T3234/T3237"
$slovakia.Range("B2").Value = "Slovakia Market"

# Remove the P32AR / P32DR rows (old rows 16 & 17) that do not apply to Slovakia.
$slovakia.Range("A16:A17").EntireRow.Delete()

# Rows 3-4 no longer need the taller wrapped height once re-entered; row 5 keeps it.
$slovakia.Rows("3:4").AutoFit()

# Restore Portugal's own selection (it is no longer the selected/active tab).
$portugal.Range("B4").Select() | Out-Null

# Make Slovakia the active (selected) sheet/tab and set its selection.
$slovakia.Activate() | Out-Null
$slovakia.Range("B4").Select() | Out-Null
